# Applies: "Updated notebook, reran simulation" -- adds two new shared strings
# ("Holden", "Rizzie Spiral"), renames "Thomas Hex" -> "Matthies Hex", reruns the
# simulation (new float values for every data row) and appends two more rows
# (Michael-CCHex, Michael-SNHex) that used to be the last two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 31,20
$data[0,0] = $null
$data[0,1] = 0
$data[0,2] = 1
$data[0,3] = 2
$data[0,4] = 3
$data[0,5] = 4
$data[0,6] = 5
$data[0,7] = 6
$data[0,8] = 7
$data[0,9] = 8
$data[0,10] = 9
$data[0,11] = 10
$data[0,12] = 11
$data[0,13] = 12
$data[0,14] = 13
$data[0,15] = 14
$data[0,16] = 15
$data[0,17] = 16
$data[0,18] = 17
$data[0,19] = 18
$data[1,0] = 0
$data[1,1] = 'HKL'
$data[1,2] = '[1, 1, 0]'
$data[1,3] = '[2, 2, 2]'
$data[1,4] = '[3, 1, 0]'
$data[1,5] = '[3, 2, 1]'
$data[1,6] = '[2, 1, 1]'
$data[1,7] = '[4, 0, 0]'
$data[1,8] = '[2, 2, 0]'
$data[1,9] = '[2, 0, 0]'
$data[1,10] = '1Pair-A'
$data[1,11] = '1Pair-B'
$data[1,12] = '2Pairs-A'
$data[1,13] = '2Pairs-B'
$data[1,14] = '3Pairs-A'
$data[1,15] = '3Pairs-B'
$data[1,16] = '3Pairs-C'
$data[1,17] = '4Pairs'
$data[1,18] = '5A4F'
$data[1,19] = 'MaxUnique'
$data[2,0] = 1
$data[2,1] = 'Spiral5'
$data[2,2] = 1.000237077428879
$data[2,3] = 1.000972743138677
$data[2,4] = 0.9992831226717243
$data[2,5] = 1.000320798069191
$data[2,6] = 1.000355253926609
$data[2,7] = 0.9987492523799008
$data[2,8] = 1.000237077428879
$data[2,9] = 0.9987492523799008
$data[2,10] = 1.000237077428879
$data[2,11] = 1.000355253926609
$data[2,12] = 0.999552253153255
$data[2,13] = 0.999552253153255
$data[2,14] = 0.9994625429927447
$data[2,15] = 0.9997805279117964
$data[2,16] = 0.9997805279117964
$data[2,17] = 0.999894665291067
$data[2,18] = 0.999894665291067
$data[2,19] = 0.999986374602497
$data[3,0] = 2
$data[3,1] = 'Holden'
$data[3,2] = 1.009516405143995
$data[3,3] = 1.04608192043703
$data[3,4] = 0.9671381326919676
$data[3,5] = 1.01465017878051
$data[3,6] = 1.016762844131858
$data[3,7] = 0.9433010373289121
$data[3,8] = 1.009516405143995
$data[3,9] = 0.9433010373289121
$data[3,10] = 1.009516405143995
$data[3,11] = 1.016762844131858
$data[3,12] = 0.980031940730385
$data[3,13] = 0.980031940730385
$data[3,14] = 0.9757340047175792
$data[3,15] = 0.9898600955349218
$data[3,16] = 0.9898600955349218
$data[3,17] = 0.9947741729371902
$data[3,18] = 0.9947741729371902
$data[3,19] = 0.9995750864190455
$data[4,0] = 3
$data[4,1] = 'Rizzie Spiral'
$data[4,2] = 1.066374473382178
$data[4,3] = 1.085827953909528
$data[4,4] = 0.9065998797135668
$data[4,5] = 1.043674698969831
$data[4,6] = 1.034333212009058
$data[4,7] = 0.8170818087847411
$data[4,8] = 1.066374473382178
$data[4,9] = 0.8170818087847411
$data[4,10] = 1.066374473382178
$data[4,11] = 1.034333212009058
$data[4,12] = 0.9257075103968998
$data[4,13] = 0.9257075103968998
$data[4,14] = 0.9193383001691222
$data[4,15] = 0.9725964980586591
$data[4,16] = 0.9725964980586591
$data[4,17] = 0.9960409918895389
$data[4,18] = 0.9960409918895389
$data[4,19] = 0.9923153377948171
$data[5,0] = 4
$data[5,1] = 'RotRing OmegaMax-90'
$data[5,2] = 1.01012504071284
$data[5,3] = 1.050594877518245
$data[5,4] = 0.9641670619023227
$data[5,5] = 1.015947703055966
$data[5,6] = 1.018343862912224
$data[5,7] = 0.9384219348019158
$data[5,8] = 1.01012504071284
$data[5,9] = 0.9384219348019158
$data[5,10] = 1.01012504071284
$data[5,11] = 1.018343862912224
$data[5,12] = 0.9783828988570701
$data[5,13] = 0.9783828988570701
$data[5,14] = 0.973644286538821
$data[5,15] = 0.9889636128089935
$data[5,16] = 0.9889636128089935
$data[5,17] = 0.9942539697849552
$data[5,18] = 0.9942539697849552
$data[5,19] = 0.9996000801505858
$data[6,0] = 5
$data[6,1] = 'Equal Angle'
$data[6,2] = 1.012768862680116
$data[6,3] = 1.054950553342937
$data[6,4] = 0.9598770629322755
$data[6,5] = 1.017945219409217
$data[6,6] = 1.020075407550432
$data[6,7] = 0.9301484846902028
$data[6,8] = 1.012768862680116
$data[6,9] = 0.9301484846902028
$data[6,10] = 1.012768862680116
$data[6,11] = 1.020075407550432
$data[6,12] = 0.9751119461203173
$data[6,13] = 0.9751119461203173
$data[6,14] = 0.9700336517243033
$data[6,15] = 0.9876642516402501
$data[6,16] = 0.9876642516402501
$data[6,17] = 0.9939404044002165
$data[6,18] = 0.9939404044002165
$data[6,19] = 0.9992942651008634
$data[7,0] = 6
$data[7,1] = 'Tilt Rotate'
$data[7,2] = 1.030993762996975
$data[7,3] = 1.192300761886641
$data[7,4] = 0.8686822421396768
$data[7,5] = 1.058159383837483
$data[7,6] = 1.069338622572605
$data[7,7] = 0.7774707567759208
$data[7,8] = 1.030993762996975
$data[7,9] = 0.7774707567759208
$data[7,10] = 1.030993762996975
$data[7,11] = 1.069338622572605
$data[7,12] = 0.9234046896742627
$data[7,13] = 0.9234046896742627
$data[7,14] = 0.9051638738294008
$data[7,15] = 0.959267714115167
$data[7,16] = 0.959267714115167
$data[7,17] = 0.9771992263356191
$data[7,18] = 0.9771992263356191
$data[7,19] = 0.9994909217015504
$data[8,0] = 7
$data[8,1] = 'CLR'
$data[8,2] = 1.002855224395957
$data[8,3] = 1.00439315918221
$data[8,4] = 0.9955772725437189
$data[8,5] = 1.002053536191132
$data[8,6] = 1.001723626506172
$data[8,7] = 0.9914951085433201
$data[8,8] = 1.002855224395957
$data[8,9] = 0.9914951085433201
$data[8,10] = 1.002855224395957
$data[8,11] = 1.001723626506172
$data[8,12] = 0.9966093675247463
$data[8,13] = 0.9966093675247463
$data[8,14] = 0.9962653358644037
$data[8,15] = 0.9986913198151499
$data[8,16] = 0.9986913198151499
$data[8,17] = 0.9997322959603517
$data[8,18] = 0.9997322959603517
$data[8,19] = 0.9996829878937518
$data[9,0] = 8
$data[9,1] = 'Rizzie Hex'
$data[9,2] = 1.000404034907337
$data[9,3] = 1.000129837303039
$data[9,4] = 0.999661090977086
$data[9,5] = 1.000165531044686
$data[9,6] = 1.000067383446498
$data[9,7] = 0.9992566704326337
$data[9,8] = 1.000404034907337
$data[9,9] = 0.9992566704326337
$data[9,10] = 1.000404034907337
$data[9,11] = 1.000067383446498
$data[9,12] = 0.9996620269395657
$data[9,13] = 0.9996620269395657
$data[9,14] = 0.9996617149520725
$data[9,15] = 0.999909362928823
$data[9,16] = 0.9999093629288228
$data[9,17] = 1.000033030923452
$data[9,18] = 1.000033030923452
$data[9,19] = 0.9999474246852132
$data[10,0] = 9
$data[10,1] = 'Matthies Hex'
$data[10,2] = 1.00461239214728
$data[10,3] = 1.007395301358965
$data[10,4] = 0.9926792649804023
$data[10,5] = 1.003394945161671
$data[10,6] = 1.002893941360726
$data[10,7] = 0.9859720275282761
$data[10,8] = 1.00461239214728
$data[10,9] = 0.9859720275282761
$data[10,10] = 1.00461239214728
$data[10,11] = 1.002893941360726
$data[10,12] = 0.9944329844445012
$data[10,13] = 0.9944329844445012
$data[10,14] = 0.9938484112898016
$data[10,15] = 0.9978261203454274
$data[10,16] = 0.9978261203454274
$data[10,17] = 0.9995226882958905
$data[10,18] = 0.9995226882958905
$data[10,19] = 0.9994913120895533
$data[11,0] = 10
$data[11,1] = 'Tilt Rotate_Partial'
$data[11,2] = 1.030980615844372
$data[11,3] = 1.195682552669861
$data[11,4] = 0.8667447591427362
$data[11,5] = 1.05899203046137
$data[11,6] = 1.070519332453391
$data[11,7] = 0.7744563423654272
$data[11,8] = 1.030980615844372
$data[11,9] = 0.7744563423654272
$data[11,10] = 1.030980615844372
$data[11,11] = 1.070519332453391
$data[11,12] = 0.9224878374094091
$data[11,13] = 0.9224878374094091
$data[11,14] = 0.9039068113205181
$data[11,15] = 0.9586520968877302
$data[11,16] = 0.9586520968877302
$data[11,17] = 0.9767342266268908
$data[11,18] = 0.9767342266268908
$data[11,19] = 0.9995626054895262
$data[12,0] = 11
$data[12,1] = 'RotRing OmegaMax-60'
$data[12,2] = 1.013325024714872
$data[12,3] = 1.05115017964073
$data[12,4] = 0.9617094649523588
$data[12,5] = 1.017180500849642
$data[12,6] = 1.018767118372519
$data[12,7] = 0.9327399195481095
$data[12,8] = 1.013325024714872
$data[12,9] = 0.9327399195481095
$data[12,10] = 1.013325024714872
$data[12,11] = 1.018767118372519
$data[12,12] = 0.9757535189603144
$data[12,13] = 0.9757535189603144
$data[12,14] = 0.9710721676243291
$data[12,15] = 0.9882773542118336
$data[12,16] = 0.9882773542118337
$data[12,17] = 0.9945392718375934
$data[12,18] = 0.9945392718375934
$data[12,19] = 0.9991453680130387
$data[13,0] = 12
$data[13,1] = 'Equal Angle_Partial'
$data[13,2] = 1.009833733494736
$data[13,3] = 1.059798864231577
$data[13,4] = 0.9590358190526305
$data[13,5] = 1.018150961894734
$data[13,6] = 1.021573688421054
$data[13,7] = 0.9304930183789486
$data[13,8] = 1.009833733494736
$data[13,9] = 0.9304930183789486
$data[13,10] = 1.009833733494736
$data[13,11] = 1.021573688421054
$data[13,12] = 0.9760333534000011
$data[13,13] = 0.9760333534000011
$data[13,14] = 0.9703675086175442
$data[13,15] = 0.9873001467649128
$data[13,16] = 0.9873001467649128
$data[13,17] = 0.9929335434473687
$data[13,18] = 0.9929335434473687
$data[13,19] = 0.9998143475789466
$data[14,0] = 13
$data[14,1] = 'Rizzie Hex_Partial'
$data[14,2] = 0.988735153240399
$data[14,3] = 0.9672212166098864
$data[14,4] = 1.026378881070176
$data[14,5] = 0.9880402441459442
$data[14,6] = 0.9877542739873784
$data[14,7] = 1.047613307877775
$data[14,8] = 0.988735153240399
$data[14,9] = 1.047613307877775
$data[14,10] = 0.988735153240399
$data[14,11] = 0.9877542739873784
$data[14,12] = 1.017683790932577
$data[14,13] = 1.017683790932577
$data[14,14] = 1.020582154311777
$data[14,15] = 1.008034245035184
$data[14,16] = 1.008034245035184
$data[14,17] = 1.003209472086488
$data[14,18] = 1.003209472086488
$data[14,19] = 1.000957179488593
$data[15,0] = 14
$data[15,1] = 'ND Single'
$data[15,2] = 1.052343800000001
$data[15,3] = 1.339376599999999
$data[15,4] = 0.7698135699999986
$data[15,5] = 1.101841100000001
$data[15,6] = 1.122210300000001
$data[15,7] = 0.6110551399999998
$data[15,8] = 1.052343800000001
$data[15,9] = 0.6110551399999998
$data[15,10] = 1.052343800000001
$data[15,11] = 1.122210300000001
$data[15,12] = 0.8666327200000004
$data[15,13] = 0.8666327200000004
$data[15,14] = 0.8343596699999999
$data[15,15] = 0.928536413333334
$data[15,16] = 0.928536413333334
$data[15,17] = 0.9594882600000009
$data[15,18] = 0.9594882600000009
$data[15,19] = 0.9994400850000001
$data[16,0] = 15
$data[16,1] = 'RD Single'
$data[16,2] = 1.1248129
$data[16,3] = 1.055856
$data[16,4] = 0.88627042
$data[16,5] = 1.0550102
$data[16,6] = 1.0262848
$data[16,7] = 0.7563714999999999
$data[16,8] = 1.1248129
$data[16,9] = 0.7563714999999999
$data[16,10] = 1.1248129
$data[16,11] = 1.0262848
$data[16,12] = 0.8913281499999999
$data[16,13] = 0.8913281499999999
$data[16,14] = 0.88964224
$data[16,15] = 0.9691563999999998
$data[16,16] = 0.9691564
$data[16,17] = 1.008070525
$data[16,18] = 1.008070525
$data[16,19] = 0.9841009700000001
$data[17,0] = 16
$data[17,1] = 'TD Single'
$data[17,2] = 1.0555695
$data[17,3] = 1.3923208
$data[17,4] = 0.73750648
$data[17,5] = 1.1157804
$data[17,6] = 1.1405585
$data[17,7] = 0.55983377
$data[17,8] = 1.0555695
$data[17,9] = 0.55983377
$data[17,10] = 1.0555695
$data[17,11] = 1.1405585
$data[17,12] = 0.850196135
$data[17,13] = 0.850196135
$data[17,14] = 0.8126329166666667
$data[17,15] = 0.9186539233333333
$data[17,16] = 0.9186539233333333
$data[17,17] = 0.9528828174999999
$data[17,18] = 0.9528828174999999
$data[17,19] = 1.000261575
$data[18,0] = 17
$data[18,1] = 'Morris Single'
$data[18,2] = 0.9791626800000001
$data[18,3] = 0.93185156
$data[18,4] = 1.0533855
$data[18,5] = 0.97579072
$data[18,6] = 0.97440308
$data[18,7] = 1.0960966
$data[18,8] = 0.9791626800000001
$data[18,9] = 1.0960966
$data[18,10] = 0.9791626800000001
$data[18,11] = 0.97440308
$data[18,12] = 1.03524984
$data[18,13] = 1.03524984
$data[18,14] = 1.04129506
$data[18,15] = 1.01655412
$data[18,16] = 1.01655412
$data[18,17] = 1.00720626
$data[18,18] = 1.00720626
$data[18,19] = 1.00178169
$data[19,0] = 18
$data[19,1] = 'Ring Perpendicular to ND'
$data[19,2] = 1.057112248493151
$data[19,3] = 1.086846673972603
$data[19,4] = 0.9120855553424658
$data[19,5] = 1.040855593424658
$data[19,6] = 1.034165629589041
$data[19,7] = 0.8306111506849313
$data[19,8] = 1.057112248493151
$data[19,9] = 0.8306111506849313
$data[19,10] = 1.057112248493151
$data[19,11] = 1.034165629589041
$data[19,12] = 0.9323883901369863
$data[19,13] = 0.9323883901369863
$data[19,14] = 0.9256207785388127
$data[19,15] = 0.9739630095890411
$data[19,16] = 0.9739630095890411
$data[19,17] = 0.9947503193150686
$data[19,18] = 0.9947503193150686
$data[19,19] = 0.993612808584475
$data[20,0] = 19
$data[20,1] = 'Ring Perpendicular to RD'
$data[20,2] = 0.9996629678947367
$data[20,3] = 1.090827807894737
$data[20,4] = 0.9482600973684209
$data[20,5] = 1.022106726842105
$data[20,6] = 1.031342835789474
$data[20,7] = 0.9204919189473684
$data[20,8] = 0.9996629678947367
$data[20,9] = 0.9204919189473684
$data[20,10] = 0.9996629678947367
$data[20,11] = 1.031342835789474
$data[20,12] = 0.9759173773684211
$data[20,13] = 0.9759173773684211
$data[20,14] = 0.9666982840350876
$data[20,15] = 0.9838325742105263
$data[20,16] = 0.9838325742105264
$data[20,17] = 0.987790172631579
$data[20,18] = 0.987790172631579
$data[20,19] = 1.00211539245614
$data[21,0] = 20
$data[21,1] = 'Ring Perpendicular to TD'
$data[21,2] = 1.066816447368421
$data[21,3] = 1.103303594736842
$data[21,4] = 0.8962658610526315
$data[21,5] = 1.04813602631579
$data[21,6] = 1.040448594736842
$data[21,7] = 0.8007687910526315
$data[21,8] = 1.066816447368421
$data[21,9] = 0.8007687910526315
$data[21,10] = 1.066816447368421
$data[21,11] = 1.040448594736842
$data[21,12] = 0.9206086928947369
$data[21,13] = 0.9206086928947369
$data[21,14] = 0.912494415614035
$data[21,15] = 0.9693446110526317
$data[21,16] = 0.9693446110526317
$data[21,17] = 0.9937125701315791
$data[21,18] = 0.9937125701315791
$data[21,19] = 0.9926232192105265
$data[22,0] = 21
$data[22,1] = 'OffsetFTD'
$data[22,2] = 0.9565558125374185
$data[22,3] = 0.8865750126537043
$data[22,4] = 1.094569509018459
$data[22,5] = 0.9568255160372351
$data[22,6] = 0.9569365107262972
$data[22,7] = 1.173467059371904
$data[22,8] = 0.9565558125374185
$data[22,9] = 1.173467059371904
$data[22,10] = 0.9565558125374185
$data[22,11] = 0.9569365107262972
$data[22,12] = 1.0652017850491
$data[22,13] = 1.0652017850491
$data[22,14] = 1.07499102637222
$data[22,15] = 1.02898646087854
$data[22,16] = 1.02898646087854
$data[22,17] = 1.010878798793259
$data[22,18] = 1.010878798793259
$data[22,19] = 1.004154903390836
$data[23,0] = 22
$data[23,1] = 'OffsetATD'
$data[23,2] = 0.9881400096302552
$data[23,3] = 0.969129600173114
$data[23,4] = 1.025175499027671
$data[23,5] = 0.988737009118141
$data[23,6] = 0.9889826844603619
$data[23,7] = 1.04444737005346
$data[23,8] = 0.9881400096302552
$data[23,9] = 1.04444737005346
$data[23,10] = 0.9881400096302552
$data[23,11] = 0.9889826844603619
$data[23,12] = 1.016715027256911
$data[23,13] = 1.016715027256911
$data[23,14] = 1.019535184513831
$data[23,15] = 1.007190021381359
$data[23,16] = 1.007190021381359
$data[23,17] = 1.002427518443583
$data[23,18] = 1.002427518443583
$data[23,19] = 1.000768695410501
$data[24,0] = 23
$data[24,1] = 'OffsetF45'
$data[24,2] = 0.9891261156197422
$data[24,3] = 0.9446912132375229
$data[24,4] = 1.038768516182305
$data[24,5] = 0.9828655059055107
$data[24,6] = 0.9802891304142625
$data[24,7] = 1.065687611343511
$data[24,8] = 0.9891261156197422
$data[24,9] = 1.065687611343511
$data[24,10] = 0.9891261156197422
$data[24,11] = 0.9802891304142625
$data[24,12] = 1.022988370878887
$data[24,13] = 1.022988370878887
$data[24,14] = 1.02824841931336
$data[24,15] = 1.011700952459172
$data[24,16] = 1.011700952459172
$data[24,17] = 1.006057243249314
$data[24,18] = 1.006057243249314
$data[24,19] = 1.000238015450476
$data[25,0] = 24
$data[25,1] = 'OffsetA45'
$data[25,2] = 0.9927659664095386
$data[25,3] = 0.991457346216166
$data[25,4] = 1.010102917857146
$data[25,5] = 0.9951083754200345
$data[25,6] = 0.9960723305446388
$data[25,7] = 1.021084726870479
$data[25,8] = 0.9927659664095386
$data[25,9] = 1.021084726870479
$data[25,10] = 0.9927659664095386
$data[25,11] = 0.9960723305446388
$data[25,12] = 1.008578528707559
$data[25,13] = 1.008578528707559
$data[25,14] = 1.009086658424088
$data[25,15] = 1.003307674608219
$data[25,16] = 1.003307674608219
$data[25,17] = 1.000672247558549
$data[25,18] = 1.000672247558549
$data[25,19] = 1.001098610553001
$data[26,0] = 25
$data[26,1] = 'OffsetFRD'
$data[26,2] = 1.05896550999756
$data[26,3] = 1.022981830125628
$data[26,4] = 0.9478335360028276
$data[26,5] = 1.025483744486642
$data[26,6] = 1.011705240838845
$data[26,7] = 0.8860492708011987
$data[26,8] = 1.05896550999756
$data[26,9] = 0.8860492708011987
$data[26,10] = 1.05896550999756
$data[26,11] = 1.011705240838845
$data[26,12] = 0.9488772558200217
$data[26,13] = 0.9488772558200217
$data[26,14] = 0.9485293492142904
$data[26,15] = 0.9855733405458676
$data[26,16] = 0.9855733405458676
$data[26,17] = 1.003921382908791
$data[26,18] = 1.003921382908791
$data[26,19] = 0.9921698553754501
$data[27,0] = 26
$data[27,1] = 'OffsetARD'
$data[27,2] = 0.9911203093065557
$data[27,3] = 1.0328201399228
$data[27,4] = 0.9867215774662738
$data[27,5] = 1.005359847239953
$data[27,6] = 1.011219747729839
$data[27,7] = 0.9833706552342959
$data[27,8] = 0.9911203093065557
$data[27,9] = 0.9833706552342959
$data[27,10] = 0.9911203093065557
$data[27,11] = 1.011219747729839
$data[27,12] = 0.9972952014820675
$data[27,13] = 0.9972952014820675
$data[27,14] = 0.9937706601434696
$data[27,15] = 0.9952369040902301
$data[27,16] = 0.9952369040902301
$data[27,17] = 0.9942077553943116
$data[27,18] = 0.9942077553943116
$data[27,19] = 1.001768712816619
$data[28,0] = 27
$data[28,1] = 'Gaussian Quadrature'
$data[28,2] = 1.007539063242566
$data[28,3] = 1.015455696747633
$data[28,4] = 0.9861260814300561
$data[28,5] = 1.00635705542387
$data[28,6] = 1.005870624744431
$data[28,7] = 0.9741932594053861
$data[28,8] = 1.007539063242566
$data[28,9] = 0.9741932594053861
$data[28,10] = 1.007539063242566
$data[28,11] = 1.005870624744431
$data[28,12] = 0.9900319420749084
$data[28,13] = 0.9900319420749084
$data[28,14] = 0.9887299885266243
$data[28,15] = 0.9958676491307944
$data[28,16] = 0.9958676491307944
$data[28,17] = 0.9987855026587373
$data[28,18] = 0.9987855026587373
$data[28,19] = 0.9992569634989903
$data[29,0] = 28
$data[29,1] = 'Michael-CCHex'
$data[29,2] = 1.006092248463665
$data[29,3] = 1.003446333253574
$data[29,4] = 0.9939092861905958
$data[29,5] = 1.002969866162964
$data[29,6] = 1.001684933311993
$data[29,7] = 0.9868604164094521
$data[29,8] = 1.006092248463665
$data[29,9] = 0.9868604164094521
$data[29,10] = 1.006092248463665
$data[29,11] = 1.001684933311993
$data[29,12] = 0.9942726748607228
$data[29,13] = 0.9942726748607228
$data[29,14] = 0.9941515453040138
$data[29,15] = 0.9982125327283704
$data[29,16] = 0.9982125327283704
$data[29,17] = 1.000182461662194
$data[29,18] = 1.000182461662194
$data[29,19] = 0.9991605139653741
$data[30,0] = 29
$data[30,1] = 'Michael-SNHex'
$data[30,2] = 0.9926498008646213
$data[30,3] = 0.9459604706558941
$data[30,4] = 1.036553947592946
$data[30,5] = 0.9836448600057938
$data[30,6] = 0.9799391249235617
$data[30,7] = 1.063023969518443
$data[30,8] = 0.9926498008646213
$data[30,9] = 1.063023969518443
$data[30,10] = 0.9926498008646213
$data[30,11] = 0.9799391249235617
$data[30,12] = 1.021481547221002
$data[30,13] = 1.021481547221002
$data[30,14] = 1.026505680678317
$data[30,15] = 1.011870965102209
$data[30,16] = 1.011870965102209
$data[30,17] = 1.007065674042812
$data[30,18] = 1.007065674042812
$data[30,19] = 1.00029536226021

$ws.Range("A1:T31").Value = $data

# New rows 30/31 need the same bold/border/centered style as the rest of column A.
$ws.Range("A29").Copy($ws.Range("A30:A31"))
$ws.Range("A30").Value = 28
$ws.Range("A31").Value = 29

Write-Output "applied"
